# Update the "Neural_Net" row (row 4, columns B:L) on each of the three
# sheets (RMSE, MAE, R2) with new values.

$wb = $excel.ActiveWorkbook

# Sheet 1: RMSE
$ws1 = $wb.Worksheets.Item("RMSE")
$ws1.Range("B4").Value = 0.000126310951
$ws1.Range("C4").Value = 0.000179093072
$ws1.Range("D4").Value = 0.0001551694364
$ws1.Range("E4").Value = 0.0001851826088
$ws1.Range("F4").Value = 0.00034677498
$ws1.Range("G4").Value = 0.000346907255
$ws1.Range("H4").Value = 0.000633107441
$ws1.Range("I4").Value = 0.000742578989
$ws1.Range("J4").Value = 0.002052861896
$ws1.Range("K4").Value = 0.001805163885
$ws1.Range("L4").Value = 0.0036330289199999

# Sheet 2: MAE
$ws2 = $wb.Worksheets.Item("MAE")
$ws2.Range("B4").Value = 0.00383735637
$ws2.Range("C4").Value = 0.00436954745
$ws2.Range("D4").Value = 0.0042803539699999
$ws2.Range("E4").Value = 0.00419651825
$ws2.Range("F4").Value = 0.0054151701
$ws2.Range("G4").Value = 0.00637234095
$ws2.Range("H4").Value = 0.00680870627
$ws2.Range("I4").Value = 0.00831050025
$ws2.Range("J4").Value = 0.01184669053
$ws2.Range("K4").Value = 0.0174808820999999
$ws2.Range("L4").Value = 0.0291258992

# Sheet 3: R2
$ws3 = $wb.Worksheets.Item("R2")
$ws3.Range("B4").Value = 0.3139212762052834
$ws3.Range("C4").Value = 0.31487076042815
$ws3.Range("D4").Value = 0.3228248654077991
$ws3.Range("E4").Value = 0.3274508537085975
$ws3.Range("F4").Value = 0.3050806206033136
$ws3.Range("G4").Value = 0.31642658225183
$ws3.Range("H4").Value = 0.2567768251908443
$ws3.Range("I4").Value = 0.2788324123804463
$ws3.Range("J4").Value = 0.1906245759994369
$ws3.Range("K4").Value = 0.0555678023953661
$ws3.Range("L4").Value = -0.3258858502682218
